$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17: round the easting/northing coordinates to whole numbers
$ws.Range("Q17").Value = 364847
$ws.Range("R17").Value = 6872339

# Rows 18-20 had their species/observation data reshuffled:
#   new row 18 <- old row 20 data
#   new row 19 <- old row 18 data
#   new row 20 <- old row 19 data
# (location/date/reporter columns for these rows stay as-is)

# New row 18 (was row 20's species data), coordinates rounded
$ws.Range("A18").Value = 112182890
$ws.Range("B18").Value = 96265
$ws.Range("D18").Value = "LC"
$ws.Range("E18").Value = 219790
$ws.Range("F18").Value = "Fläcknycklar"
$ws.Range("G18").Value = "Dactylorhiza maculata"
$ws.Range("H18").Value = "(L.) Soó"
$ws.Range("Q18").Value = 364947
$ws.Range("R18").Value = 6872308

# New row 19 (was row 18's species data), coordinates rounded
$ws.Range("A19").Value = 112181755
$ws.Range("B19").Value = 78579
$ws.Range("D19").Value = "NT"
$ws.Range("E19").Value = 2081
$ws.Range("F19").Value = "Skrovellav"
$ws.Range("G19").Value = "Lobaria scrobiculata"
$ws.Range("H19").Value = "(Scop.) DC."
$ws.Range("Q19").Value = 364894
$ws.Range("R19").Value = 6872300

# New row 20 (was row 19's species data), coordinates rounded
$ws.Range("A20").Value = 112182654
$ws.Range("B20").Value = 76918
$ws.Range("D20").Value = "NT"
$ws.Range("E20").Value = 6437
$ws.Range("F20").Value = "Blanksvart spiklav"
$ws.Range("G20").Value = "Calicium denigratum"
$ws.Range("H20").Value = "(Vain.) Tibell"
$ws.Range("Q20").Value = 364914
$ws.Range("R20").Value = 6872133

# Row 21: round the easting/northing coordinates to whole numbers
$ws.Range("Q21").Value = 364898
$ws.Range("R21").Value = 6872201

# Row 22: round the easting/northing coordinates to whole numbers
$ws.Range("Q22").Value = 364938
$ws.Range("R22").Value = 6872236

# Rows 18-22 no longer carry a start/end time (Z/AB columns) value
$ws.Range("Z18").ClearContents()
$ws.Range("AB18").ClearContents()
$ws.Range("Z19").ClearContents()
$ws.Range("AB19").ClearContents()
$ws.Range("Z20").ClearContents()
$ws.Range("AB20").ClearContents()
$ws.Range("Z21").ClearContents()
$ws.Range("AB21").ClearContents()
$ws.Range("Z22").ClearContents()
$ws.Range("AB22").ClearContents()
